$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F
$ws.Cells.Item(1, 6).Value = "time_taken"

# Copy header style (bold, centered, bordered) from E1 so F1 matches the other headers
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# time_taken values for rows 2-72 (one per data row)
$timeTaken = @(
  "2021-10-05 10:49:56.710560",
  "2021-10-05 10:49:56.710570",
  "2021-10-05 10:49:56.710574",
  "2021-10-05 10:49:56.710576",
  "2021-10-05 10:49:56.710579",
  "2021-10-05 10:49:56.710582",
  "2021-10-05 10:49:56.710584",
  "2021-10-05 10:49:56.710587",
  "2021-10-05 10:49:56.710590",
  "2021-10-05 10:49:56.710592",
  "2021-10-05 10:49:56.710595",
  "2021-10-05 10:49:56.710597",
  "2021-10-05 10:49:56.710599",
  "2021-10-05 10:49:56.710602",
  "2021-10-05 10:49:56.710604",
  "2021-10-05 10:49:56.710607",
  "2021-10-05 10:49:56.710609",
  "2021-10-05 10:49:56.710612",
  "2021-10-05 10:49:56.710614",
  "2021-10-05 10:49:56.710617",
  "2021-10-05 10:49:56.710620",
  "2021-10-05 10:49:56.710622",
  "2021-10-05 10:49:56.710624",
  "2021-10-05 10:49:56.710627",
  "2021-10-05 10:49:56.710630",
  "2021-10-05 10:49:56.710632",
  "2021-10-05 10:49:56.710635",
  "2021-10-05 10:49:56.710637",
  "2021-10-05 10:49:56.710640",
  "2021-10-05 10:49:56.710643",
  "2021-10-05 10:49:56.710645",
  "2021-10-05 10:49:56.710648",
  "2021-10-05 10:49:56.710650",
  "2021-10-05 10:49:56.710653",
  "2021-10-05 10:49:56.710655",
  "2021-10-05 10:49:56.710658",
  "2021-10-05 10:49:56.710660",
  "2021-10-05 10:49:56.710663",
  "2021-10-05 10:49:56.710665",
  "2021-10-05 10:49:56.710668",
  "2021-10-05 10:49:56.710670",
  "2021-10-05 10:49:56.710673",
  "2021-10-05 10:49:56.710675",
  "2021-10-05 10:49:56.710678",
  "2021-10-05 10:49:56.710680",
  "2021-10-05 10:49:56.710683",
  "2021-10-05 10:49:56.710685",
  "2021-10-05 10:49:56.710687",
  "2021-10-05 10:49:56.710690",
  "2021-10-05 10:49:56.710692",
  "2021-10-05 10:49:56.710695",
  "2021-10-05 10:49:56.710697",
  "2021-10-05 10:49:56.710700",
  "2021-10-05 10:49:56.710703",
  "2021-10-05 10:49:56.710705",
  "2021-10-05 10:49:56.710708",
  "2021-10-05 10:49:56.710710",
  "2021-10-05 10:49:56.710712",
  "2021-10-05 10:49:56.710715",
  "2021-10-05 10:49:56.710717",
  "2021-10-05 10:49:56.710720",
  "2021-10-05 10:49:56.710722",
  "2021-10-05 10:49:56.710725",
  "2021-10-05 10:49:56.710727",
  "2021-10-05 10:49:56.710730",
  "2021-10-05 10:49:56.710733",
  "2021-10-05 10:49:56.710736",
  "2021-10-05 10:49:56.710738",
  "2021-10-05 10:49:56.710741",
  "2021-10-05 10:49:56.710743",
  "2021-10-05 10:49:56.710745"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
